# Kumiko Example.xlsx — record piece measurements for the first two pieces.
# (Corresponds to the GUI now writing Time/Material readings to the sheet
#  whenever a new piece is logged, alongside the angle/position/length
#  validation + red-flash error textfield added in the app.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Piece #1 (row 2): Time = 1.217, Material = -1
$ws.Cells.Item(2, 3).Value = 1.2170000076293945
$ws.Cells.Item(2, 4).Value = -1.0

# Piece #2 (row 3): Time = 2.556, Material = -1
$ws.Cells.Item(3, 3).Value = 2.555999994277954
$ws.Cells.Item(3, 4).Value = -1.0
